# Duplicate_Transactions.xlsx
# "Added analytics and move valiator duplicate to All records data"
#
# 1) Fix B73 so the phone number is stored as a number instead of text.
# 2) Append the newly identified duplicate rows (74-84) to the bottom of
#    the "All records" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    # Force a literal string even when the text looks like a number or a
    # date (e.g. phone numbers, "YYYY-MM-DD" strings) while leaving the
    # cell's style untouched.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# --- 1. Correct the phone number in B73 so it is numeric, not text ---
$ws.Cells.Item(73, 2).Value = 12096508105

# --- 2. Append the newly analysed duplicate rows (74-83) ---
# Columns: A Name, B Phone, C Address, D Book, E Language,
#          K Duplicate_Reason, L Attempt_Date, M Campaign_Date, N Status

$rows = @(
    @{ Row=74; Name="Madhukar Verma"; Phone=2065044242; Address="42729 Mayfair Park Ave Fremont Fremont 94538 California USA"; Book="JKR"; Language="English"; Attempt="2025-09-15 16:36:51" },
    @{ Row=75; Name="Fnu Balan"; Phone=2065044242; Address="202 Hovis Rd Hovis Rd 28164 Nc USA"; Book="JKR"; Language="English"; Attempt="2025-09-15 16:36:57" },
    @{ Row=76; Name="Amanda Father - Vazquez"; Phone=2065044242; Address="116 cypress Vallejo Ca 94590 United States"; Book="GG"; Language="English"; Attempt="2025-09-15 16:36:59" },
    @{ Row=77; Name="Thomas Schenck"; Phone=2065044242; Address="408 210 5th avenue south Saint Petersburg 33701 Florida United States"; Book="GG"; Language="English"; Attempt="2025-09-15 16:37:01" },
    @{ Row=78; Name="Atha Bass"; Phone=2065044242; Address="202 10404 Salvia Street, Charlotte, 28277, North Carolina, USA"; Book="JKR"; Language="English"; Attempt="2025-09-15 16:37:08" },
    @{ Row=79; Name="Thyagarajan Iyer"; Phone=2065044242; Address="6092 Elmbridge Dr, San Jose, 95129, CA, USA"; Book="JKR"; Language="English"; Attempt="2025-09-15 16:37:10" },
    @{ Row=80; Name="Janakkumar Babulal"; Phone=2065044242; Address="3220 1st Street, Rosenberg, Fort Bend, TX 77471, USA"; Book="NAN"; Language="Gujarati"; Attempt="2025-09-15 16:37:12" },
    @{ Row=81; Name="Janakkumar Babulal"; Phone=9165478955; Address="3220 1st Street, Rosenberg, Fort Bend, TX 77471, USA"; Book="NAN"; Language="Gujarati"; Attempt="2025-09-15 16:37:14" },
    @{ Row=82; Name="Jatinder Das"; Phone=12096891489; Address="2027 Westmora ave.`nStockton, CA 95210"; Book="YBB"; Language="Nan"; Attempt="2025-09-15 16:37:16" },
    @{ Row=83; Name="Gurinder Singh Aulakh"; Phone=12096508105; Address="698 N Plumas Dr, Mountain House, CA 95391"; Book="JKR"; Language="Punjabi"; Attempt="2025-09-15 16:37:18" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Name
    $ws.Cells.Item($row, 2).Value = $r.Phone
    $ws.Cells.Item($row, 3).Value = $r.Address
    $ws.Cells.Item($row, 4).Value = $r.Book
    $ws.Cells.Item($row, 5).Value = $r.Language
    $ws.Cells.Item($row, 11).Value = "Same book already sent"
    $ws.Cells.Item($row, 12).Value = $r.Attempt
    Set-TextValue $ws.Cells.Item($row, 13) "2025-09-15"
    $ws.Cells.Item($row, 14).Value = "Blocked"
}

# Row 84 is a brand new (non-duplicate-of-existing-row) record; its phone
# number stays as text, matching the source export.
$ws.Cells.Item(84, 1).Value = "Sandipkumar Bhupendralal Kapadia"
Set-TextValue $ws.Cells.Item(84, 2) "2063260971"
$ws.Cells.Item(84, 3).Value = "605 West Third Street, Donalsonville, Georgia, 39845, USA"
$ws.Cells.Item(84, 4).Value = "NAN"
$ws.Cells.Item(84, 5).Value = "Gujarati"
$ws.Cells.Item(84, 11).Value = "Same book already sent"
$ws.Cells.Item(84, 12).Value = "2025-09-15 16:37:21"
Set-TextValue $ws.Cells.Item(84, 13) "2025-09-15"
$ws.Cells.Item(84, 14).Value = "Blocked"
